$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Southwest Asian / N. African Ancestry',
    'Afghan',
    'Algerian',
    'Amazigh',
    'Arab',
    'Arabic',
    'Armenian',
    'Assyrian',
    'Azerbaijan',
    'Bahraini',
    'Berber',
    'Chaldean',
    'Copt',
    'Djiboutian',
    'Egyptian',
    'Emirati',
    'Iranian',
    'Iraqi',
    'Israeli',
    'Jordanian',
    'Kurdish',
    'Kuwaiti',
    'Lebanese',
    'Levantine',
    'Libyan',
    'Mauritanian',
    'Mid East',
    'Mideast',
    'Middle Eastern',
    'Middle Eastern or North African Responses, not elsewhere classified',
    'Moroccan',
    'North African',
    'Omani',
    'Other Arab',
    'Palestinian',
    'Qatari',
    'Saudi',
    'Saudi Arabian',
    'Somali',
    'Southwest Asian or North African',
    'Sudanese',
    'SWANA',
    'Syriac',
    'Syrian',
    'Tunisian',
    'Turkish',
    'West Asian',
    'Yazidi',
    'Yemeni'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
